$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; C="1001"; D="2218686.79"},
    @{Row=4; C="1049"; D="3760636.47"},
    @{Row=6; C="691"; D="2283574.78"},
    @{Row=7; C="18"; D="39498.41"},
    @{Row=8; C="37"; D="164144.45"},
    @{Row=9; C="197"; D="585215.60"},
    @{Row=10; C="384"; D="1449720.18"},
    @{Row=11; C="177"; D="658476.34"},
    @{Row=31; C="589"; D="2496919.89"},
    @{Row=41; C="203"; D="562055.23"},
    @{Row=42; C="95"; D="440699.98"},
    @{Row=43; C="141"; D="580972.25"},
    @{Row=46; C="419"; D="1160354.43"},
    @{Row=48; C="657"; D="2714879.99"},
    @{Row=49; C="449"; D="1617996.50"},
    @{Row=52; C="3828"; D="8824141.47"},
    @{Row=57; C="4195"; D="13339755.27"},
    @{Row=63; C="132"; D="291450.00"},
    @{Row=64; C="223"; D="607326.60"},
    @{Row=65; C="160"; D="396518.88"},
    @{Row=77; C="951"; D="3351064.26"},
    @{Row=78; C="534"; D="1769407.38"},
    @{Row=87; C="229"; D="539969.00"},
    @{Row=89; C="524"; D="1889604.72"},
    @{Row=90; C="191"; D="572880.81"}
)

foreach ($u in $updates) {
    $cCell = $ws.Cells.Item($u.Row, 3)
    $cStyle = $cCell.Style
    $cCell.NumberFormat = "@"
    $cCell.Value = $u.C
    $cCell.Style = $cStyle

    $dCell = $ws.Cells.Item($u.Row, 4)
    $dStyle = $dCell.Style
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $dCell.Style = $dStyle
}
